# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 11.25 = 45674.43 pesos`n✅ 45674.43 pesos = 11.19 = 973.58 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $text

# --- tasas: update the rate cells N10/O10/N12/O12 ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 88.89
$ws2.Range("O10").Value = 4060
$ws2.Range("N12").Value = 4081.5
$ws2.Range("O12").Value = 87
